$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.672.35'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.05%  '
$ws.Range("D3").Value = '''1.591.01'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.70%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''211.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.50%  '
$ws.Range("D6").Value = '''0.510'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.24%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '''0.249'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.72%  '
$ws.Range("E9").Value = '  -1.65%  '
$ws.Range("D10").Value = '''19.69'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.60%  '
$ws.Range("E11").Value = '  -1.71%  '
$ws.Range("D12").Value = '''1.813.20'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.68%  '
$ws.Range("D13").Value = '''1.589.57'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.92%  '
$ws.Range("E14").Value = '  -2.66%  '
$ws.Range("D15").Value = '''0.527'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.62%  '
$ws.Range("D16").Value = '''64.55'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.24%  '
$ws.Range("D17").Value = '''26.648.37'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.02%  '
$ws.Range("D18").Value = '''0.0₃0728'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.86%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '''207.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.89%  '
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").Value = '''1.00'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("D21").Value = '''6.74'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.38%  '
$ws.Range("D22").Value = '''4.25'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.41%  '
$ws.Range("E23").Value = '  -2.38%  '
$ws.Range("D24").Value = '''8.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.24%  '
$ws.Range("D25").Value = '''146.96'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.51%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("E28").Value = '  -3.69%  '
$ws.Range("D29").Value = '''15.28'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.28%  '
$ws.Range("D30").Value = '''0.0504'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.50%  '
$ws.Range("E31").Value = '  -2.02%  '
$ws.Range("E32").Value = '  -4.63%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '''0.652'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +19.06%  '
$ws.Range("B34").Value = 'Maker'
$ws.Range("C34").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D34").Value = '''1.321.81'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").Value = '''2.90'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.33%  '
$ws.Range("E36").Value = '  -4.62%  '
$ws.Range("E37").Value = '  -1.86%  '
$ws.Range("E38").Value = '  -2.14%  '
$ws.Range("D39").Value = '''0.829'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.57%  '
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '''0.790'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.65%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '''5.38'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.70%  '
$ws.Range("E43").Value = '  -3.58%  '
$ws.Range("D44").Value = '''63.48'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.46%  '
$ws.Range("D45").Value = '''1.726.46'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.47%  '
$ws.Range("D46").Value = '''90.02'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.18%  '
$ws.Range("E47").Value = '  -1.34%  '
$ws.Range("D48").Value = '''0.837'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.49%  '
$ws.Range("D49").Value = '''0.0513'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.33%  '
$ws.Range("D50").Value = '''0.0973'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.69%  '
$ws.Range("E51").Value = '  -1.07%  '
